$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first "00010012 / 100" line item is removed; the remaining rows
# (00010012/500, 90351051/100, 90351051/1000) shift up to take its place,
# which also shrinks the used range from A1:B5 down to A1:B4.
$ws.Range("A2:B2").EntireRow.Delete()
